$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.712993860244751
$ws.Cells.Item(2, 5).Value = 236.9274748885455
$ws.Cells.Item(2, 6).Value = 0.007686309101501964
$ws.Cells.Item(2, 7).Value = 0.006533694960710819
$ws.Cells.Item(2, 8).Value = 0.006094620780313737
$ws.Cells.Item(2, 9).Value = 0.005842421531522889
$ws.Cells.Item(2, 10).Value = 0.005703255215933124
$ws.Cells.Item(2, 11).Value = 0.005703255215933124
$ws.Cells.Item(2, 12).Value = 0.005542350352134235
$ws.Cells.Item(2, 13).Value = 0.005248600080509847
$ws.Cells.Item(2, 14).Value = 0.005248600080509847
$ws.Cells.Item(2, 15).Value = 0.005173508915236756
$ws.Cells.Item(2, 16).Value = 0.004956574545403267
$ws.Cells.Item(2, 17).Value = 0.004913951033372017
$ws.Cells.Item(2, 18).Value = 0.004913951033372017
$ws.Cells.Item(2, 19).Value = 0.004810199110640724
$ws.Cells.Item(2, 20).Value = 0.004786028739330623
$ws.Cells.Item(2, 21).Value = 0.004722752451528996
$ws.Cells.Item(2, 22).Value = 0.004667559021229846
$ws.Cells.Item(2, 23).Value = 0.004667559021229846
$ws.Cells.Item(2, 24).Value = 0.004641101421244707
$ws.Cells.Item(2, 25).Value = 0.004618469296073011

$ws.Cells.Item(3, 3).Value = 0.817997932434082
$ws.Cells.Item(3, 5).Value = 238.4336981396937
$ws.Cells.Item(3, 6).Value = 0.007837855268614205
$ws.Cells.Item(3, 7).Value = 0.006520280231572867
$ws.Cells.Item(3, 8).Value = 0.006271053087771808
$ws.Cells.Item(3, 9).Value = 0.006221731626953616
$ws.Cells.Item(3, 10).Value = 0.006221731626953616
$ws.Cells.Item(3, 11).Value = 0.006016489405433195
$ws.Cells.Item(3, 12).Value = 0.005632985933106465
$ws.Cells.Item(3, 13).Value = 0.005632985933106465
$ws.Cells.Item(3, 14).Value = 0.005350353788437475
$ws.Cells.Item(3, 15).Value = 0.005045898160385714
$ws.Cells.Item(3, 16).Value = 0.004960643367975983
$ws.Cells.Item(3, 17).Value = 0.004960643367975983
$ws.Cells.Item(3, 18).Value = 0.004905566266833805
$ws.Cells.Item(3, 19).Value = 0.004869671144612178
$ws.Cells.Item(3, 20).Value = 0.004869671144612178
$ws.Cells.Item(3, 21).Value = 0.004766953827636977
$ws.Cells.Item(3, 22).Value = 0.004732969213924328
$ws.Cells.Item(3, 23).Value = 0.004706590131929728
$ws.Cells.Item(3, 24).Value = 0.004677876315039696
$ws.Cells.Item(3, 25).Value = 0.004647830373093443

$ws.Cells.Item(4, 3).Value = 0.6820440292358398
$ws.Cells.Item(4, 5).Value = 234.2631290752615
$ws.Cells.Item(4, 6).Value = 0.007456591039384232
$ws.Cells.Item(4, 7).Value = 0.006004470023957867
$ws.Cells.Item(4, 8).Value = 0.005496283750764196
$ws.Cells.Item(4, 9).Value = 0.005465568915706755
$ws.Cells.Item(4, 10).Value = 0.005465568915706755
$ws.Cells.Item(4, 11).Value = 0.005360917787056003
$ws.Cells.Item(4, 12).Value = 0.005343444677131004
$ws.Cells.Item(4, 13).Value = 0.005343444677131004
$ws.Cells.Item(4, 14).Value = 0.005222897906695431
$ws.Cells.Item(4, 15).Value = 0.005159688344847345
$ws.Cells.Item(4, 16).Value = 0.004865230170873073
$ws.Cells.Item(4, 17).Value = 0.004865230170873073
$ws.Cells.Item(4, 18).Value = 0.004865230170873073
$ws.Cells.Item(4, 19).Value = 0.004865230170873073
$ws.Cells.Item(4, 20).Value = 0.00485506442450627
$ws.Cells.Item(4, 21).Value = 0.004757570827791565
$ws.Cells.Item(4, 22).Value = 0.004666630887258591
$ws.Cells.Item(4, 23).Value = 0.00457697635107241
$ws.Cells.Item(4, 24).Value = 0.00457697635107241
$ws.Cells.Item(4, 25).Value = 0.00456653273051192

$ws.Cells.Item(5, 3).Value = 0.8549585342407227
$ws.Cells.Item(5, 5).Value = 222.8647225860896
$ws.Cells.Item(5, 6).Value = 0.007730477571098601
$ws.Cells.Item(5, 7).Value = 0.006395203333345155
$ws.Cells.Item(5, 8).Value = 0.005985591471392316
$ws.Cells.Item(5, 9).Value = 0.005656059168542877
$ws.Cells.Item(5, 10).Value = 0.005451352868619392
$ws.Cells.Item(5, 11).Value = 0.005451352868619392
$ws.Cells.Item(5, 12).Value = 0.005298089167732846
$ws.Cells.Item(5, 13).Value = 0.005095977851119578
$ws.Cells.Item(5, 14).Value = 0.00497801417960015
$ws.Cells.Item(5, 15).Value = 0.00497801417960015
$ws.Cells.Item(5, 16).Value = 0.00497801417960015
$ws.Cells.Item(5, 17).Value = 0.004795695806911205
$ws.Cells.Item(5, 18).Value = 0.004663805833663192
$ws.Cells.Item(5, 19).Value = 0.00444179977383707
$ws.Cells.Item(5, 20).Value = 0.00444179977383707
$ws.Cells.Item(5, 21).Value = 0.00444179977383707
$ws.Cells.Item(5, 22).Value = 0.00444179977383707
$ws.Cells.Item(5, 23).Value = 0.004436788898094684
$ws.Cells.Item(5, 24).Value = 0.004399496394142412
$ws.Cells.Item(5, 25).Value = 0.004344341570878938

$ws.Cells.Item(6, 3).Value = 1.011003017425537
$ws.Cells.Item(6, 5).Value = 236.690439047723
$ws.Cells.Item(6, 6).Value = 0.008166719753212833
$ws.Cells.Item(6, 7).Value = 0.006519477631430377
$ws.Cells.Item(6, 8).Value = 0.005872449866674677
$ws.Cells.Item(6, 9).Value = 0.005872449866674677
$ws.Cells.Item(6, 10).Value = 0.005872449866674677
$ws.Cells.Item(6, 11).Value = 0.005382240935539396
$ws.Cells.Item(6, 12).Value = 0.005382240935539396
$ws.Cells.Item(6, 13).Value = 0.005382240935539396
$ws.Cells.Item(6, 14).Value = 0.005146682192068277
$ws.Cells.Item(6, 15).Value = 0.005146682192068277
$ws.Cells.Item(6, 16).Value = 0.005071762901822975
$ws.Cells.Item(6, 17).Value = 0.004930238474648538
$ws.Cells.Item(6, 18).Value = 0.004930238474648538
$ws.Cells.Item(6, 19).Value = 0.004651014200129541
$ws.Cells.Item(6, 20).Value = 0.004651014200129541
$ws.Cells.Item(6, 21).Value = 0.004651014200129541
$ws.Cells.Item(6, 22).Value = 0.004645401434464309
$ws.Cells.Item(6, 23).Value = 0.004633104895405529
$ws.Cells.Item(6, 24).Value = 0.004627768006951929
$ws.Cells.Item(6, 25).Value = 0.004613848714380564

$ws.Cells.Item(7, 3).Value = 0.8429965972900391
$ws.Cells.Item(7, 5).Value = 239.2141299962277
$ws.Cells.Item(7, 6).Value = 0.007868996010751408
$ws.Cells.Item(7, 7).Value = 0.006431766394324435
$ws.Cells.Item(7, 8).Value = 0.005792221777710348
$ws.Cells.Item(7, 9).Value = 0.005792221777710348
$ws.Cells.Item(7, 10).Value = 0.00563403316077411
$ws.Cells.Item(7, 11).Value = 0.00563403316077411
$ws.Cells.Item(7, 12).Value = 0.00556972764686344
$ws.Cells.Item(7, 13).Value = 0.005486134467479614
$ws.Cells.Item(7, 14).Value = 0.005304254764232644
$ws.Cells.Item(7, 15).Value = 0.005304254764232644
$ws.Cells.Item(7, 16).Value = 0.00521354498223236
$ws.Cells.Item(7, 17).Value = 0.005027122390913159
$ws.Cells.Item(7, 18).Value = 0.004856852543837597
$ws.Cells.Item(7, 19).Value = 0.004856852543837597
$ws.Cells.Item(7, 20).Value = 0.004856852543837597
$ws.Cells.Item(7, 21).Value = 0.004797871257888564
$ws.Cells.Item(7, 22).Value = 0.004740458992612961
$ws.Cells.Item(7, 23).Value = 0.004697940020248485
$ws.Cells.Item(7, 24).Value = 0.004682197839501753
$ws.Cells.Item(7, 25).Value = 0.00466304346971204

$ws.Cells.Item(8, 3).Value = 0.760040283203125
$ws.Cells.Item(8, 5).Value = 233.8686316967996
$ws.Cells.Item(8, 6).Value = 0.007514702547482954
$ws.Cells.Item(8, 7).Value = 0.006287850202304441
$ws.Cells.Item(8, 8).Value = 0.005778960528254524
$ws.Cells.Item(8, 9).Value = 0.005748341007915723
$ws.Cells.Item(8, 10).Value = 0.005718128859134806
$ws.Cells.Item(8, 11).Value = 0.005718128859134806
$ws.Cells.Item(8, 12).Value = 0.005694590436215702
$ws.Cells.Item(8, 13).Value = 0.005530432345126436
$ws.Cells.Item(8, 14).Value = 0.00517468703067612
$ws.Cells.Item(8, 15).Value = 0.005096379075953043
$ws.Cells.Item(8, 16).Value = 0.005027437251443174
$ws.Cells.Item(8, 17).Value = 0.004982737987644092
$ws.Cells.Item(8, 18).Value = 0.004844333999460456
$ws.Cells.Item(8, 19).Value = 0.00478569835377601
$ws.Cells.Item(8, 20).Value = 0.004774407115374098
$ws.Cells.Item(8, 21).Value = 0.004663452789606281
$ws.Cells.Item(8, 22).Value = 0.004663452789606281
$ws.Cells.Item(8, 23).Value = 0.004624395055431951
$ws.Cells.Item(8, 24).Value = 0.004606195217820604
$ws.Cells.Item(8, 25).Value = 0.004558842723134495

$ws.Cells.Item(9, 3).Value = 0.8079986572265625
$ws.Cells.Item(9, 5).Value = 246.6821658874032
$ws.Cells.Item(9, 6).Value = 0.008202759255531005
$ws.Cells.Item(9, 7).Value = 0.006473146220843957
$ws.Cells.Item(9, 8).Value = 0.006031423245740067
$ws.Cells.Item(9, 9).Value = 0.005899763504592122
$ws.Cells.Item(9, 10).Value = 0.005771505267934555
$ws.Cells.Item(9, 11).Value = 0.00539415244125093
$ws.Cells.Item(9, 12).Value = 0.00539415244125093
$ws.Cells.Item(9, 13).Value = 0.00539415244125093
$ws.Cells.Item(9, 14).Value = 0.00539415244125093
$ws.Cells.Item(9, 15).Value = 0.00539415244125093
$ws.Cells.Item(9, 16).Value = 0.005364492553468776
$ws.Cells.Item(9, 17).Value = 0.005329619605195014
$ws.Cells.Item(9, 18).Value = 0.005195719072856938
$ws.Cells.Item(9, 19).Value = 0.005009423986870507
$ws.Cells.Item(9, 20).Value = 0.004997231847881914
$ws.Cells.Item(9, 21).Value = 0.004997231847881914
$ws.Cells.Item(9, 22).Value = 0.004860337663604265
$ws.Cells.Item(9, 23).Value = 0.004860337663604265
$ws.Cells.Item(9, 24).Value = 0.004833934164909617
$ws.Cells.Item(9, 25).Value = 0.004808619218078035

$ws.Cells.Item(10, 3).Value = 0.8010013103485107
$ws.Cells.Item(10, 5).Value = 250.2937931478282
$ws.Cells.Item(10, 6).Value = 0.006773381411383712
$ws.Cells.Item(10, 7).Value = 0.006382137391800865
$ws.Cells.Item(10, 8).Value = 0.006376698607663674
$ws.Cells.Item(10, 9).Value = 0.00635710619006945
$ws.Cells.Item(10, 10).Value = 0.005747669227463939
$ws.Cells.Item(10, 11).Value = 0.005747669227463939
$ws.Cells.Item(10, 12).Value = 0.00561224571417935
$ws.Cells.Item(10, 13).Value = 0.00561224571417935
$ws.Cells.Item(10, 14).Value = 0.005325555254018083
$ws.Cells.Item(10, 15).Value = 0.005325555254018083
$ws.Cells.Item(10, 16).Value = 0.005297166536433169
$ws.Cells.Item(10, 17).Value = 0.005262974217904259
$ws.Cells.Item(10, 18).Value = 0.00514217869409893
$ws.Cells.Item(10, 19).Value = 0.00508248608118054
$ws.Cells.Item(10, 20).Value = 0.004986343885243593
$ws.Cells.Item(10, 21).Value = 0.004986343885243593
$ws.Cells.Item(10, 22).Value = 0.004967006790082992
$ws.Cells.Item(10, 23).Value = 0.004916718966646104
$ws.Cells.Item(10, 24).Value = 0.004879021308924525
$ws.Cells.Item(10, 25).Value = 0.004879021308924525

$ws.Cells.Item(11, 3).Value = 0.811039924621582
$ws.Cells.Item(11, 5).Value = 243.9252197103633
$ws.Cells.Item(11, 6).Value = 0.007370076423741731
$ws.Cells.Item(11, 7).Value = 0.006063030851256914
$ws.Cells.Item(11, 8).Value = 0.006063030851256914
$ws.Cells.Item(11, 9).Value = 0.005821899190410497
$ws.Cells.Item(11, 10).Value = 0.005821899190410497
$ws.Cells.Item(11, 11).Value = 0.005348601249616374
$ws.Cells.Item(11, 12).Value = 0.005348601249616374
$ws.Cells.Item(11, 13).Value = 0.005348601249616374
$ws.Cells.Item(11, 14).Value = 0.005262465484090115
$ws.Cells.Item(11, 15).Value = 0.005230424521522148
$ws.Cells.Item(11, 16).Value = 0.005230424521522148
$ws.Cells.Item(11, 17).Value = 0.005060253196097795
$ws.Cells.Item(11, 18).Value = 0.005040117774048682
$ws.Cells.Item(11, 19).Value = 0.004880733215122886
$ws.Cells.Item(11, 20).Value = 0.004880733215122886
$ws.Cells.Item(11, 21).Value = 0.004880733215122886
$ws.Cells.Item(11, 22).Value = 0.004859336524651299
$ws.Cells.Item(11, 23).Value = 0.004797941993235856
$ws.Cells.Item(11, 24).Value = 0.004797941993235856
$ws.Cells.Item(11, 25).Value = 0.004754877577200063
